# Fix NaN values in machine cost tables by correcting the Economic Life
# (and its dependent computed values) for the Fecon FTX128 and Kaiser S2-2
# machines from 10 years down to 5 years.

$wb = $excel.ActiveWorkbook

# --- Sheet "Summary Costs" ---
$wsSummary = $wb.Worksheets.Item("Summary Costs")
$wsSummary.Range("E2").Value = 5
$wsSummary.Range("I2").Value = 70.8685879997
$wsSummary.Range("E3").Value = 5
$wsSummary.Range("I3").Value = 114.191413638

# --- Sheet "Operating Costs" ---
$wsOperating = $wb.Worksheets.Item("Operating Costs")
$wsOperating.Range("C2").Value = 28.962611901
$wsOperating.Range("C3").Value = 54.1600842549

# --- Sheet "Fixed Costs" ---
$wsFixed = $wb.Worksheets.Item("Fixed Costs")
$wsFixed.Range("C2").Value = 32000
$wsFixed.Range("D2").Value = 16320
$wsFixed.Range("E2").Value = 4080
$wsFixed.Range("F2").Value = 4080

$wsFixed.Range("C3").Value = 59840
$wsFixed.Range("D3").Value = 30518.4
$wsFixed.Range("E3").Value = 7629.6
$wsFixed.Range("F3").Value = 7629.6
